$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''42.812.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -4.91%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.220.73'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -6.17%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.03%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''316.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +1.66%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''99.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -7.49%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.593'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -5.93%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.02%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  -7.78%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''37.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -8.67%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''53.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -3.16%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0829'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -9.43%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''7.82'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -7.13%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  -3.31%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '''2.557.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -6.16%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.861'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -11.32%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  -6.14%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''2.218.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -6.21%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').Value = '''15.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +8.04%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '''42.768.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -5.00%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''0.0₃0964'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -8.79%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''6.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -10.77%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -10.70%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''3.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -8.90%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''236.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -8.67%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  -7.62%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.15%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''10.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -9.02%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  -5.01%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -11.42%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  -8.06%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  -7.77%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''34.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -7.84%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''156.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -7.02%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -6.22%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''3.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +9.72%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''1.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +13.64%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.123'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -5.44%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''4.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -6.30%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''3.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.14%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  -11.40%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  -7.67%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''1.920.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +2.22%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  +0.06%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''12.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -3.90%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''89.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -10.72%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  -9.01%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''5.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -3.72%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''60.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -12.69%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''75.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -7.08%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.861'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +17.57%  '
$ws.Range('E51').Style = 'Normal'
